# Applies the "statement_32.xlsx" bank-statement content swap described in the diff:
# new account holder / card number, new statement period dates, new transactions
# for rows 6-8, row 9's transaction is removed (blanked out), and the running
# balances / next-billing-date footer are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block: account holder name + card number
$ws.Range("C2").Value = "Hartmut"
# Card number looks numeric but must stay text; a leading apostrophe forces
# Excel to store it as text (quote-prefixed) instead of coercing to a number.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 05.10.2023"

# Transaction row 6
$ws.Range("B6").Value = "07.10."
$ws.Range("C6").Value = "08.10."
$ws.Range("D6").Value = "MCDONALDS Viechtach"
$ws.Range("E6").Value = "15,95-"

# Transaction row 7
$ws.Range("B7").Value = "09.10."
$ws.Range("C7").Value = "10.10."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "59,34-"

# Transaction row 8
$ws.Range("B8").Value = "10.10."
$ws.Range("C8").Value = "11.10."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 44232469"
$ws.Range("E8").Value = "41,67-"

# Row 9's transaction is dropped entirely in the new statement - blank it out
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
# Match the blank-row formatting used for the amount column elsewhere (center/center/wrap)
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Closing balance + next statement date footer
$ws.Range("D12").Value = "KONTOSTAND AM 12.10.2023"
$ws.Range("E12").Value = "116,96-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 17.10.2023"
